$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ42858302"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1181.624098109003
$ws.Range("B3").Value = -181.5941779923623
$ws.Range("B4").Value = 208.8379114156564
$ws.Range("B5").Value = -0.5491321936600122
$ws.Range("B6").Value = 0.0001871643399941286
$ws.Range("B7").Value = 134.0122331843442
$ws.Range("B8").Value = -0.01272236054498421
$ws.Range("B9").Value = -3291.680634991939
$ws.Range("B10").Value = 845.6802874912155
$ws.Range("B11").Value = 812.4016832777647
$ws.Range("B12").Value = -15.54461918017296

$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ43034502"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 4237.001053598709
$ws.Range("B3").Value = -204.0649470628337
$ws.Range("B4").Value = 744.9300690634615
$ws.Range("B5").Value = -0.2191497692295745
$ws.Range("B6").Value = 0.0004824707445811862
$ws.Range("B7").Value = 132.6025589026348
$ws.Range("B8").Value = -19.82738406880108
$ws.Range("B9").Value = -4870.965326972306
$ws.Range("B10").Value = -3234.509111248539
$ws.Range("B11").Value = -3673.745869926614
$ws.Range("B12").Value = -71.43382755899302

$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ43362110"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1848.060313508171
$ws.Range("B3").Value = -595.6755285506224
$ws.Range("B4").Value = 93.23117051693498
$ws.Range("B5").Value = -0.6011660935781107
$ws.Range("B6").Value = 0.00005641032158212986
$ws.Range("B7").Value = 77.60210194698911
$ws.Range("B8").Value = -11.67119927371259
$ws.Range("B9").Value = -320.6346891488647
$ws.Range("B10").Value = -1434.366859384931
$ws.Range("B11").Value = -7233.572630725179
$ws.Range("B12").Value = 92.63484971576754

$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ43558219"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 10069.54645508002
$ws.Range("B3").Value = 320.1622053487449
$ws.Range("B4").Value = 425.9972729528349
$ws.Range("B5").Value = -0.5801736267249433
$ws.Range("B6").Value = 0.0001185184699125768
$ws.Range("B7").Value = 3.422859933417953
$ws.Range("B8").Value = -42.88186444845755
$ws.Range("B9").Value = 2548.916807831914
$ws.Range("B10").Value = 4399.099947331744
$ws.Range("B11").Value = 1548.348753884994
$ws.Range("B12").Value = -55.93471031170793

$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ43720325"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 3820.26103938583
$ws.Range("B3").Value = 881.7831535256573
$ws.Range("B4").Value = 958.0326139439544
$ws.Range("B5").Value = -0.1606599029481099
$ws.Range("B6").Value = 0.0003465338622374232
$ws.Range("B7").Value = -9.910030039282674
$ws.Range("B8").Value = -88.00234810359245
$ws.Range("B9").Value = -1445.718969340433
$ws.Range("B10").Value = -1311.705479630753
$ws.Range("B11").Value = -699.1246305258865
$ws.Range("B12").Value = 125.5603289813146

$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ43890944"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 9153.500300985688
$ws.Range("B3").Value = -64.14357371301617
$ws.Range("B4").Value = 304.7762982578702
$ws.Range("B5").Value = -0.3340749518522051
$ws.Range("B6").Value = 0.0002637036648642402
$ws.Range("B7").Value = 137.9608803820523
$ws.Range("B8").Value = 10.58182470058178
$ws.Range("B9").Value = -6223.885782365833
$ws.Range("B10").Value = -10468.91531591885
$ws.Range("B11").Value = -4527.976423054981
$ws.Range("B12").Value = -158.093138886856

$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ44065315"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1295.085764128828
$ws.Range("B3").Value = -15.1416961803493
$ws.Range("B4").Value = 431.8006325276965
$ws.Range("B5").Value = -0.2135882921902499
$ws.Range("B6").Value = 0.0001816012568217496
$ws.Range("B7").Value = 61.70142797060316
$ws.Range("B8").Value = -16.91476796652108
$ws.Range("B9").Value = -3545.791314386668
$ws.Range("B10").Value = 1729.326817111401
$ws.Range("B11").Value = 1209.943565147319
$ws.Range("B12").Value = 28.87259173308541

$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ44231412"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1707.717043474387
$ws.Range("B3").Value = -457.1655015047227
$ws.Range("B4").Value = 493.2658944900413
$ws.Range("B5").Value = -0.3358739310308021
$ws.Range("B6").Value = 0.0001163202741252781
$ws.Range("B7").Value = 100.7833037927397
$ws.Range("B8").Value = -22.1953458914731
$ws.Range("B9").Value = -1620.514741445678
$ws.Range("B10").Value = -1064.471442833834
$ws.Range("B11").Value = 308.7711155969755
$ws.Range("B12").Value = 38.39057682756595

$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ44395309"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 14842.32038908326
$ws.Range("B3").Value = -698.6257530603559
$ws.Range("B4").Value = 72.57242624404506
$ws.Range("B5").Value = -0.3321180927277555
$ws.Range("B6").Value = -0.000326072295703561
$ws.Range("B7").Value = -137.1903292037812
$ws.Range("B8").Value = -33.91327910128518
$ws.Range("B9").Value = -3444.30802720323
$ws.Range("B10").Value = 10277.62146368493
$ws.Range("B11").Value = 12017.80378990486
$ws.Range("B12").Value = -46.74181036243755

$ws = $wb.Worksheets.Item(10)
$ws.Name = "summ44562648"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = -3702.407448003949
$ws.Range("B3").Value = -150.1310491276381
$ws.Range("B4").Value = 1671.343234454274
$ws.Range("B5").Value = -0.6544000076915073
$ws.Range("B6").Value = 0.0006698861522281352
$ws.Range("B7").Value = 298.9930486271155
$ws.Range("B8").Value = -16.21873099964709
$ws.Range("B9").Value = -3976.087848362391
$ws.Range("B10").Value = 1994.050869021958
$ws.Range("B11").Value = -834.9859691595623
$ws.Range("B12").Value = -122.2847163798979

$ws = $wb.Worksheets.Item(11)
$ws.Name = "summ44733656"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2606.371353364451
$ws.Range("B3").Value = 1242.036563052839
$ws.Range("B4").Value = -80.41324947612111
$ws.Range("B5").Value = -1.948686313519683
$ws.Range("B6").Value = 0.0003394768351965532
$ws.Range("B7").Value = 114.6753482580017
$ws.Range("B8").Value = 19.17128191916979
$ws.Range("B9").Value = 4738.965068281464
$ws.Range("B10").Value = 2035.363662833191
$ws.Range("B11").Value = 1805.541819851196
$ws.Range("B12").Value = -66.12917916868236

$ws = $wb.Worksheets.Item(12)
$ws.Name = "summ44904594"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = -1504.847226627149
$ws.Range("B3").Value = 1593.412910738571
$ws.Range("B4").Value = 1284.071800912127
$ws.Range("B5").Value = 1.109675878321736
$ws.Range("B6").Value = 0.0003736937959880069
$ws.Range("B7").Value = -33.60259526285085
$ws.Range("B8").Value = -90.94016986811687
$ws.Range("B9").Value = -8703.023151590489
$ws.Range("B10").Value = -4744.508197203941
$ws.Range("B11").Value = -14662.85026892589
$ws.Range("B12").Value = 209.7677355875867

$ws = $wb.Worksheets.Item(13)
$ws.Name = "summ45081570"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 3046.477136688158
$ws.Range("B3").Value = -254.0511750718492
$ws.Range("B4").Value = 613.7443970462723
$ws.Range("B5").Value = -0.1475130570324947
$ws.Range("B6").Value = 0.0001417775343751783
$ws.Range("B7").Value = 89.75542415515713
$ws.Range("B8").Value = -32.14569973458015
$ws.Range("B9").Value = -603.492843144777
$ws.Range("B10").Value = -1721.887097523721
$ws.Range("B11").Value = -536.0667673316107
$ws.Range("B12").Value = 20.92107833026273

$ws = $wb.Worksheets.Item(14)
$ws.Name = "summ45253951"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 4920.744704131179
$ws.Range("B3").Value = -181.3313144965314
$ws.Range("B4").Value = 2016.033643749437
$ws.Range("B5").Value = -1.193564781986594
$ws.Range("B6").Value = 0.001645887655883472
$ws.Range("B7").Value = 429.9609183633106
$ws.Range("B8").Value = -1.829682350009847
$ws.Range("B9").Value = -8387.670924508744
$ws.Range("B10").Value = -5967.410153033499
$ws.Range("B11").Value = -5955.422179239427
$ws.Range("B12").Value = -437.4997734557978

$ws = $wb.Worksheets.Item(15)
$ws.Name = "summ45428063"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1392.559874865631
$ws.Range("B3").Value = 240.0187672624215
$ws.Range("B4").Value = 287.3143217764028
$ws.Range("B5").Value = -0.3029791336091865
$ws.Range("B6").Value = 0.000167193136363937
$ws.Range("B7").Value = 80.4121972046529
$ws.Range("B8").Value = -4.25666992991053
$ws.Range("B9").Value = -2727.610935188491
$ws.Range("B10").Value = -1120.728244098307
$ws.Range("B11").Value = 1784.942339876608
$ws.Range("B12").Value = -2.578914646367934

$ws = $wb.Worksheets.Item(16)
$ws.Name = "summ45602576"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 6082.018787707406
$ws.Range("B3").Value = -337.9576888574984
$ws.Range("B4").Value = 419.8087165784536
$ws.Range("B5").Value = -0.655452571846276
$ws.Range("B6").Value = 0.000382839501157035
$ws.Range("B7").Value = 133.5149160380013
$ws.Range("B8").Value = -37.56688907216937
$ws.Range("B9").Value = -3045.230338791183
$ws.Range("B10").Value = 3422.964830330355
$ws.Range("B11").Value = -7704.289964695963
$ws.Range("B12").Value = -18.18500608947602

$ws = $wb.Worksheets.Item(17)
$ws.Name = "summ45780097"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 9628.22841275266
$ws.Range("B3").Value = 171.0066324880295
$ws.Range("B4").Value = 199.1292810488795
$ws.Range("B5").Value = -0.6969725453540399
$ws.Range("B6").Value = 0.00008102036132450038
$ws.Range("B7").Value = 76.61479529520304
$ws.Range("B8").Value = -5.57072337195045
$ws.Range("B9").Value = -5245.467026563352
$ws.Range("B10").Value = -1546.516450315188
$ws.Range("B11").Value = 2149.282051796955
$ws.Range("B12").Value = -100.97298463222

$ws = $wb.Worksheets.Item(18)
$ws.Name = "summ45964106"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2592.549490690606
$ws.Range("B3").Value = -67.62381019227269
$ws.Range("B4").Value = 347.8206299235808
$ws.Range("B5").Value = 0.1241721049558162
$ws.Range("B6").Value = 0.00004827566035152148
$ws.Range("B7").Value = 38.17356934460497
$ws.Range("B8").Value = -8.778591281692286
$ws.Range("B9").Value = -4361.063374782148
$ws.Range("B10").Value = -818.2187390786603
$ws.Range("B11").Value = 3638.801447790203
$ws.Range("B12").Value = 0.0280687172868852

$ws = $wb.Worksheets.Item(19)
$ws.Name = "summ46162215"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2313.849322625108
$ws.Range("B3").Value = -161.1678437198416
$ws.Range("B4").Value = 407.1930991802301
$ws.Range("B5").Value = 0.4517583455933778
$ws.Range("B6").Value = -0.0001241899582371461
$ws.Range("B7").Value = -15.29104117990278
$ws.Range("B8").Value = 10.27908624367392
$ws.Range("B9").Value = -2483.902262674801
$ws.Range("B10").Value = -6143.779626830677
$ws.Range("B11").Value = -2131.956067184816
$ws.Range("B12").Value = 20.77398114724667

$ws = $wb.Worksheets.Item(20)
$ws.Name = "summ46344136"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = -322.3595467149062
$ws.Range("B3").Value = -702.0290672777082
$ws.Range("B4").Value = 486.9802572498463
$ws.Range("B5").Value = -0.242346982718489
$ws.Range("B6").Value = 0.0002292638558567227
$ws.Range("B7").Value = 79.66838356552637
$ws.Range("B8").Value = -16.98428140902124
$ws.Range("B9").Value = -5124.722763588215
$ws.Range("B10").Value = 5781.699603151414
$ws.Range("B11").Value = 2825.194358999619
$ws.Range("B12").Value = 52.99230605047251

$ws = $wb.Worksheets.Item(21)
$ws.Name = "summ46559080"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 8646.839333103584
$ws.Range("B3").Value = -162.7184311414971
$ws.Range("B4").Value = 476.081581411122
$ws.Range("B5").Value = -0.3196502845786036
$ws.Range("B6").Value = 0.00004443289343568206
$ws.Range("B7").Value = 33.38594020133644
$ws.Range("B8").Value = -26.98796858287513
$ws.Range("B9").Value = -1748.570393299864
$ws.Range("B10").Value = 1596.883795514048
$ws.Range("B11").Value = -3903.899186213632
$ws.Range("B12").Value = -42.87149846197997

$ws = $wb.Worksheets.Item(22)
$ws.Name = "summ46729074"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2373.929539742311
$ws.Range("B3").Value = -249.2901625379229
$ws.Range("B4").Value = 559.6586589488661
$ws.Range("B5").Value = -0.1700644288163691
$ws.Range("B6").Value = 0.0001420681409736079
$ws.Range("B7").Value = 91.04627582984514
$ws.Range("B8").Value = -28.7429928866527
$ws.Range("B9").Value = -642.5809135754935
$ws.Range("B10").Value = -1409.927350941576
$ws.Range("B11").Value = 1092.206011065615
$ws.Range("B12").Value = 22.99250124822362

$ws = $wb.Worksheets.Item(23)
$ws.Name = "summ46901548"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2225.996930316884
$ws.Range("B3").Value = 215.7672883661162
$ws.Range("B4").Value = 349.5780705375936
$ws.Range("B5").Value = 0.6882332908572959
$ws.Range("B6").Value = -0.0005815066116593229
$ws.Range("B7").Value = -197.8865860486517
$ws.Range("B8").Value = 25.6110730414764
$ws.Range("B9").Value = 2235.13360886545
$ws.Range("B10").Value = -1072.593221392975
$ws.Range("B11").Value = 757.6421589351776
$ws.Range("B12").Value = 76.78682892883012

$ws = $wb.Worksheets.Item(24)
$ws.Name = "summ47098530"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 2592.813494419323
$ws.Range("B3").Value = 176.7146370975206
$ws.Range("B4").Value = 1367.786461138123
$ws.Range("B5").Value = -0.7754800825002892
$ws.Range("B6").Value = 0.0005188845923609652
$ws.Range("B7").Value = 230.6893170845594
$ws.Range("B8").Value = 1.28107450303429
$ws.Range("B9").Value = -1087.109810396418
$ws.Range("B10").Value = 2159.814126787962
$ws.Range("B11").Value = -22931.59822176286
$ws.Range("B12").Value = -169.6219269208668

$ws = $wb.Worksheets.Item(25)
$ws.Name = "summ47303944"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = -1847.837090759611
$ws.Range("B3").Value = -871.5275246819784
$ws.Range("B4").Value = 515.8102885725202
$ws.Range("B5").Value = -1.577954279213698
$ws.Range("B6").Value = 0.0006001318326188266
$ws.Range("B7").Value = 267.1696070361258
$ws.Range("B8").Value = -48.31844708511889
$ws.Range("B9").Value = -2361.074397457942
$ws.Range("B10").Value = 5799.166778666564
$ws.Range("B11").Value = 644.4985987583424
$ws.Range("B12").Value = 82.99623672179214

$ws = $wb.Worksheets.Item(26)
$ws.Name = "summ47491424"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 206.7461525454655
$ws.Range("B3").Value = -417.369173350756
$ws.Range("B4").Value = 155.0429429827946
$ws.Range("B5").Value = -0.08539778848684199
$ws.Range("B6").Value = -0.00006467069149705316
$ws.Range("B7").Value = 52.74086235761381
$ws.Range("B8").Value = -1.666040671612166
$ws.Range("B9").Value = -3056.245743370593
$ws.Range("B10").Value = 872.2588547876153
$ws.Range("B11").Value = 88.13017942770375
$ws.Range("B12").Value = 70.18699223875194

$ws = $wb.Worksheets.Item(27)
$ws.Name = "summ47666656"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 1192.695161910076
$ws.Range("B3").Value = -267.987046733424
$ws.Range("B4").Value = 677.4349960035952
$ws.Range("B5").Value = 0.1830937128888599
$ws.Range("B6").Value = 0.0001479709518502114
$ws.Range("B7").Value = 94.24219369429616
$ws.Range("B8").Value = -31.62323249395173
$ws.Range("B9").Value = -2224.798734771349
$ws.Range("B10").Value = -3711.178368597211
$ws.Range("B11").Value = -456.7647562281873
$ws.Range("B12").Value = 42.34952913043523

$ws = $wb.Worksheets.Item(28)
$ws.Name = "summ47838393"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 3422.248101969792
$ws.Range("B3").Value = -93.76410051753555
$ws.Range("B4").Value = 942.6073394600041
$ws.Range("B5").Value = -0.01427395453988134
$ws.Range("B6").Value = 0.0002196307547319004
$ws.Range("B7").Value = 110.3950093025624
$ws.Range("B8").Value = -52.40216268519033
$ws.Range("B9").Value = 571.77222908767
$ws.Range("B10").Value = -2320.372698309042
$ws.Range("B11").Value = 1013.04445095044
$ws.Range("B12").Value = 6.490173628079788

$ws = $wb.Worksheets.Item(29)
$ws.Name = "summ48003824"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 12469.86190381955
$ws.Range("B3").Value = -788.2369516602948
$ws.Range("B4").Value = -1969.311594636489
$ws.Range("B5").Value = -8.037873176188732
$ws.Range("B6").Value = -0.001238597852768636
$ws.Range("B7").Value = -582.8046946988234
$ws.Range("B8").Value = -293.4495076341514
$ws.Range("B9").Value = 35153.38444250217
$ws.Range("B10").Value = 86360.37920834392
$ws.Range("B11").Value = 80450.91541081818
$ws.Range("B12").Value = 992.987663779485

$ws = $wb.Worksheets.Item(30)
$ws.Name = "summ48176359"
$ws.Rows.Item(11).Delete()
$ws.Range("B2").Value = 687.6624341271418
$ws.Range("B3").Value = -658.4013083053906
$ws.Range("B4").Value = 571.6146753643453
$ws.Range("B5").Value = -0.439300938286098
$ws.Range("B6").Value = 0.0004756077414100018
$ws.Range("B7").Value = 258.0271179347651
$ws.Range("B8").Value = -49.00150483294783
$ws.Range("B9").Value = -7860.73846937636
$ws.Range("B10").Value = 1321.57125893519
$ws.Range("B11").Value = -13471.69085308108
$ws.Range("B12").Value = 49.14517521947641
